$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10, column C: N/A -> Andre, Kaj
$ws.Range("C10").Value = "Andre, Kaj "

# Clear the two template rows (11-12) so their placeholder strings are fully freed
$ws.Range("A11:G12").Value = ""

# Row 11: new diary entry (Jan 16th)
$ws.Range("A11").Value = "Jan 16th"
$ws.Range("B11").Value = "5pm-8pm"
$ws.Range("C11").Value = "Andre, Kaj, Ping"
$ws.Range("D11").Value = "Learn how to read code from academic view and industry view."
$ws.Range("E11").Value = "Understand different way to debug an exist project. "
$ws.Range("F11").Value = "1. Error happened when building jpacman2 and solved it by redownloading .`n2. For debuging in a project, class names and method names will be useful."
$ws.Range("G11").Value = "Cool"

# Row 12: new diary entry (Jan 16th, 4pm-5pm)
$ws.Range("A12").Value = "Jan 16th"
$ws.Range("B12").Value = "4pm-5pm"
$ws.Range("C12").Value = "Junxian, Wenchia"
$ws.Range("D12").Value = "Build JEdit on Intellij."
$ws.Range("E12").Value = "Build Jedit, and studey Ant, Maven, "
$ws.Range("F12").Value = "Error happened when building Jedit using Ant as some dependencies are missing."
$ws.Range("G12").Value = "Feeling difficult to understand why to use Ant Maven and their difference."

# Row heights
$ws.Range("A11").RowHeight = 102
$ws.Range("A12").RowHeight = 51

# Selection
$ws.Range("B11").Select()
